$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New time-log rows (49-56) -------------------------------------------------
# Each row: Task/Activity, Doc/Component, 'Check in', When (date), Who, How Long
# Column D (date) re-uses the existing date style (copied from D44, numFmt
# "m/d/yyyy") so no stray custom number-format gets created.

function Add-LogRow {
    param($Row, $A, $B, $C, $Year, $Month, $Day, $E, $F)

    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C

    $ws.Range("D44").Copy()
    $ws.Cells.Item($Row, 4).PasteSpecial(-4122)
    $ws.Cells.Item($Row, 4).Value = (Get-Date -Year $Year -Month $Month -Day $Day -Hour 0 -Minute 0 -Second 0)

    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
}

Add-LogRow 49 "analysis and solving readable tree format" "CID 8c5b78b" "Java code;" 2013 12 12 "Susan" 2

Add-LogRow 50 "functional test runs and captures for test plan scripts" "CID 4ea9f4f" "FunctionalTestRunsFitness.xlsx" 2013 12 13 "Susan" 3
$ws.Rows(50).RowHeight = 25.5

Add-LogRow 51 "analysis of usable print test results in fitness evaluator spreadsheets" "CID b0f3d7f" "CodeTestingEarlyAnalysis.xlsx" 2013 12 12 "Susan" 2
$ws.Rows(51).RowHeight = 25.5

Add-LogRow 52 "check in mutation update from Li" "CID 95d0b09" "mutation.java" 2013 12 13 "Susan" 0.1
$ws.Rows(52).RowHeight = 12.75

Add-LogRow 53 "final design doc & add to report" "CID 37c8bc1" "Final Class Diagram visio, report" 2013 12 13 "Susan" 1.5
$ws.Rows(53).RowHeight = 12.75

Add-LogRow 54 "code merge correction" "CID d61fbe1" "Code check-in collision" 2013 12 13 "Susan" 1
$ws.Rows(54).RowHeight = 12.75

Add-LogRow 55 "test run: functional testing to analyse bad fitness values" "CID 72e7c66" "FuncitnalTestRunsFitness.xlsx" 2013 12 13 "Susan" 1
$ws.Rows(55).RowHeight = 12.75

Add-LogRow 56 "analysis and review of the fitness output defect" "NA" "Ongoing - Li and Susan together" 2013 12 13 "Li & Susan" 1
$ws.Rows(56).RowHeight = 12.75

# --- View state: scroll position + active selection ---------------------------
$ws.Range("A45").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
